$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.51%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.301"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.58%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08377"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.60%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.944"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.36%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9722"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.95%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.07%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1127"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.68%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1906"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.52%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09673"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.34%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04568"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.31%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1062"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.03%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.98%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005776"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.56%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.366"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.02%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.427"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.08%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3357"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.16%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.387"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-18.53%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1351"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.34%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04174"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.68%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001239"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.62%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004429"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.12%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001302"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.77%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002983"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.19%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02706"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-4.52%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05630"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.38%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007801"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.12%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.31%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007313"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.89%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002123"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.64%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007907"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.40%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3508"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006903"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.87%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.20%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003494"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.32%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003535"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.08%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.20%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.20%"
